$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wording in C43 ("one the" -> "one of the")
$ws.Range("C43").Value = 'But he could well have become one of the 9 million children under the age of 5 who die each year, mostly from preventable and treatable afflictions.'

# Prepare new rows 45-57: copy date-style from B43 and text-style from C43
$ws.Range("B43").Copy()
$ws.Range("B45:B57").PasteSpecial(-4122)
$ws.Range("C43").Copy()
$ws.Range("C45:C57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set B-column dates first (all rows)
$ws.Range("B45").Value2 = 43360
$ws.Range("B46").Value2 = 43360
$ws.Range("B47").Value2 = 43360
$ws.Range("B48").Value2 = 43360
$ws.Range("B49").Value2 = 43360
$ws.Range("B50").Value2 = 43360
$ws.Range("B51").Value2 = 43360
$ws.Range("B52").Value2 = 43360
$ws.Range("B53").Value2 = 43360
$ws.Range("B54").Value2 = 43360
$ws.Range("B55").Value2 = 43360
$ws.Range("B56").Value2 = 43360
$ws.Range("B57").Value2 = 43360

# Set C-column text in the exact order the new strings were authored (shared-string append order)
$ws.Range("C47").Value = 'So we must shape the world that he deserves.'
$ws.Range("C48").Value = 'Half of humanity lives on less than 2.5 dollars a day.'
$ws.Range("C46").Value = 'That little boy''s future is tied to ours; Our security is ultimately linked to his well-being.'
$ws.Range("C49").Value = 'That child deserves a world without extreme hunger and dependence that it fosters.'
$ws.Range("C51").Value = 'Yet Africa''s crop production remains the lowest in the world. With your generation''s leadership and ingenuity, you can make it the highest.'
$ws.Range("C50").Value = 'Agricultural research has produced stronger crops that yield more, adapt faster, and better resist drought, disease, and pests.'
$ws.Range("C52").Value = 'a quality education'
$ws.Range("C53").Value = 'human trafficking'
$ws.Range("C55").Value = 'new vaccines for tuberculosis'
$ws.Range("C56").Value = 'smart therapies that kill cancer cells and leave their healthy neighbors untouched'
$ws.Range("C54").Value = 'new cures for old plagues'
$ws.Range("C57").Value = 'needle-free immunizations to stop pandemics in their tracks.'
$ws.Range("C45").Value = 'They are both children of God, of equal worth, equal consequence, and equal rights.'

$ws.Range("C45:C57").Select()
